# Add a "Save" column (H) to the s_vals sheet, mirroring the header
# style of the existing "sum" column (G) and filling in the two data
# rows with the Save flag values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy G1's formatting (bold font, border, centered alignment) into H1
# so the new header cell matches the style of the other header cells.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data values for the Save column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
